# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns.
# D/E cells hold plain text in the workbook (prices use "." as a thousands
# separator so they are not valid numbers anyway; percentages are padded
# with spaces), so values that DO look numeric ("0.614", "59.39", ...) are
# written with a leading apostrophe to force Excel to keep them as text,
# matching the original inline-string storage instead of silently turning
# them into numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.791.46"
$ws.Range("E2").Value = "  -0.29%  "

$ws.Range("D3").Value = "2.028.88"
$ws.Range("E3").Value = "  -1.03%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("E5").Value = "  -0.92%  "

$ws.Range("D6").Value = "'0.614"
$ws.Range("E6").Value = "  -0.30%  "

$ws.Range("D7").Value = "'59.39"
$ws.Range("E7").Value = "  +2.40%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "'0.384"
$ws.Range("E9").Value = "  -0.40%  "

$ws.Range("D10").Value = "'0.0811"
$ws.Range("E10").Value = "  +0.70%  "

$ws.Range("D11").Value = "'0.104"
$ws.Range("E11").Value = "  +0.39%  "

$ws.Range("D12").Value = "'14.57"
$ws.Range("E12").Value = "  +0.35%  "

$ws.Range("D13").Value = "2.330.00"
$ws.Range("E13").Value = "  -0.94%  "

$ws.Range("D14").Value = "'21.09"
$ws.Range("E14").Value = "  +2.43%  "

$ws.Range("D15").Value = "'0.759"
$ws.Range("E15").Value = "  +1.17%  "

$ws.Range("D16").Value = "'5.17"
$ws.Range("E16").Value = "  -1.63%  "

$ws.Range("D17").Value = "2.041.34"
$ws.Range("E17").Value = "  -0.90%  "

$ws.Range("D18").Value = "37.749.52"
$ws.Range("E18").Value = "  -0.19%  "

$ws.Range("D19").Value = "'6.01"
$ws.Range("E19").Value = "  -1.75%  "

$ws.Range("D20").Value = "'69.93"
$ws.Range("E20").Value = "  +0.46%  "

$ws.Range("D21").Value = "0.0₃0824"
$ws.Range("E21").Value = "  -0.66%  "

$ws.Range("D22").Value = "'224.93"
$ws.Range("E22").Value = "  +0.20%  "

$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("D24").Value = "'2.38"
$ws.Range("E24").Value = "  -2.81%  "

$ws.Range("E25").Value = "  -1.43%  "

$ws.Range("D26").Value = "'9.27"
$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("D27").Value = "'165.15"
$ws.Range("E27").Value = "  -0.64%  "

$ws.Range("E28").Value = "  -2.22%  "

$ws.Range("E29").Value = "  -0.19%  "

$ws.Range("E30").Value = "  -4.60%  "

$ws.Range("E31").Value = "  +0.90%  "

$ws.Range("D32").Value = "'4.43"
$ws.Range("E32").Value = "  -1.95%  "

$ws.Range("D33").Value = "'2.09"
$ws.Range("E33").Value = "  +2.07%  "

$ws.Range("E34").Value = "  -1.31%  "

$ws.Range("E35").Value = "  -1.42%  "

$ws.Range("E36").Value = "  +7.23%  "

$ws.Range("D37").Value = "'2.25"
$ws.Range("E37").Value = "  -3.39%  "

$ws.Range("D38").Value = "'3.24"
$ws.Range("E38").Value = "  -1.86%  "

$ws.Range("E39").Value = "  +0.07%  "

$ws.Range("D40").Value = "1.521.09"
$ws.Range("E40").Value = "  +2.53%  "

$ws.Range("E41").Value = "  +0.44%  "

$ws.Range("D42").Value = "'96.62"
$ws.Range("E42").Value = "  -0.88%  "

$ws.Range("D43").Value = "'16.75"
$ws.Range("E43").Value = "  +0.62%  "

$ws.Range("E44").Value = "  -0.67%  "

$ws.Range("D45").Value = "'0.0917"
$ws.Range("E45").Value = "  -1.58%  "

$ws.Range("E46").Value = "  -1.47%  "

$ws.Range("D47").Value = "'4.06"
$ws.Range("E47").Value = "  -5.47%  "

$ws.Range("E48").Value = "  -0.81%  "

$ws.Range("D49").Value = "'2.95"
$ws.Range("E49").Value = "  -0.20%  "

$ws.Range("D50").Value = "'7.06"
$ws.Range("E50").Value = "  +1.04%  "

$ws.Range("D51").Value = "2.218.69"
$ws.Range("E51").Value = "  -1.10%  "
